$wb = $excel.ActiveWorkbook

# Update both the "展览" sheet and its duplicate "全部类型" sheet:
#   F3 (想去人数 for row "丽水·龙泉ACG动漫游戏博览会"): 1667 -> 1670
#   F6 (想去人数 for row "丽水·第四届HP国风动漫游戏嘉年华"): 442 -> 447
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1670
    $ws.Range("F6").Value = 447
}
